# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.333.66"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.645.86"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.51"
$ws.Range("E5").Value = "  -1.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.24"
$ws.Range("E6").Value = "  -3.03%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.644.77"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("E11").Value = "  +1.35%  "

$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.96"
$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.132.32"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("E15").Value = "  -2.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.214.31"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("E17").Value = "  -2.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.661.82"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("E19").Value = "  +2.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.01"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.09"
$ws.Range("E21").Value = "  -2.26%  "

$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.47"
$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  -3.03%  "

$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.783.55"
$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0953"
$ws.Range("E30").Value = "  -0.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.06"
$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "496.39"
$ws.Range("E32").Value = "  -5.28%  "

$ws.Range("E33").Value = "  -2.58%  "

$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.87"
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("E37").Value = "  +3.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.36"
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  -6.93%  "

$ws.Range("E43").Value = "  -3.39%  "

$ws.Range("E44").Value = "  -3.71%  "

$ws.Range("E45").Value = "  -1.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.17"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.94"
$ws.Range("E47").Value = "  +1.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.551"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.67"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("E51").Value = "  -1.11%  "
